$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.354.24'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -3.76%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.852.22'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -5.23%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.55%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '322.58'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.80%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.46%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4478'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -5.99%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3824'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -5.32%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.18'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -9.83%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07860'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -6.55%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.013'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -4.02%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.37'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -4.20%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.866.12'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -4.02%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.853'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -4.75%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.116'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -6.07%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.72%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '85.93'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -5.55%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001030'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -3.70%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06505'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.01'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -8.10%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.41%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.462'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -6.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.343.54'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -3.84%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.76'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -6.19%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.253'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.49%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.068.42'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -4.52%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '151.65'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.18%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.37'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.98%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.059'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -4.44%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.521'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -6.76%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '119.62'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.18%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09319'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.99%  '
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9350'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -4.55%  '
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.465'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.98%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.603'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.57%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.247'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -6.29%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02222'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -4.51%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05961'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -4.20%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.203'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.58%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.277'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -6.97%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.48%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5896'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -4.86%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -3.71%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -7.91%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.258'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -6.90%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5648'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -5.05%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.22'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -5.88%  '
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.925'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -6.46%  '
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'PancakeSwap'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.362'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.79%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06858'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.77%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '108.22'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.07%  '
